$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "02/16/2021"
$ws.Range("J2").Value = 80
$ws.Range("M2").Value = 10375
$ws.Range("O2").Value = 'Provincia de Quillota'

# Row 3
$ws.Range("D3").Value = "04/29/2021"
$ws.Range("J3").Value = 25
$ws.Range("K3").Value = 10000
$ws.Range("L3").Value = 10000
$ws.Range("M3").Value = 10000
$ws.Range("P3").Value = 167

# Row 4
$ws.Range("D4").Value = "05/17/2021"
$ws.Range("J4").Value = 25
$ws.Range("L4").Value = 11000
$ws.Range("M4").Value = 10400
$ws.Range("P4").Value = 173

# Row 5
$ws.Range("D5").Value = "04/05/2021"
$ws.Range("J5").Value = 20
$ws.Range("K5").Value = 9000
$ws.Range("L5").Value = 9000
$ws.Range("M5").Value = 9000
$ws.Range("P5").Value = 150

# Row 6
$ws.Range("D6").Value = "01/04/2021"
$ws.Range("J6").Value = 10
$ws.Range("K6").Value = 9000
$ws.Range("L6").Value = 9000
$ws.Range("M6").Value = 9000
$ws.Range("P6").Value = 150

# Row 7
$ws.Range("D7").Value = "03/29/2021"
$ws.Range("J7").Value = 35

# Row 8
$ws.Range("D8").Value = "12/14/2020"
$ws.Range("J8").Value = 15
$ws.Range("K8").Value = 7000
$ws.Range("L8").Value = 7000
$ws.Range("M8").Value = 7000
$ws.Range("O8").Value = 'Provincia de Limarí'
$ws.Range("P8").Value = 117

# Row 9
$ws.Range("D9").Value = "12/21/2020"
$ws.Range("J9").Value = 15
$ws.Range("K9").Value = 7000
$ws.Range("L9").Value = 7000
$ws.Range("M9").Value = 7000
$ws.Range("P9").Value = 117

# Row 10
$ws.Range("D10").Value = "07/28/2021"
$ws.Range("J10").Value = 45
$ws.Range("N10").Value = '$/caja 50 unidades'
$ws.Range("O10").Value = 'Provincia de Quillota'
$ws.Range("P10").Value = 180
$ws.Range("Q10").Value = 50

# Row 11
$ws.Range("D11").Value = "03/22/2021"
$ws.Range("J11").Value = 25
$ws.Range("K11").Value = 10000
$ws.Range("L11").Value = 10000
$ws.Range("M11").Value = 10000
$ws.Range("N11").Value = '$/caja 60 unidades'
$ws.Range("O11").Value = 'Provincia de Limarí'
$ws.Range("P11").Value = 167
$ws.Range("Q11").Value = 60

# Row 12
$ws.Range("D12").Value = "04/26/2021"
$ws.Range("J12").Value = 30
$ws.Range("K12").Value = 10000
$ws.Range("L12").Value = 10000
$ws.Range("M12").Value = 10000
$ws.Range("P12").Value = 167

